$wb = $excel.ActiveWorkbook

# --- Services sheet: move selection off the Services tab before we
# activate FirewallPolicies, so its sheetView loses tabSelected and
# lands on the new selected cell (A2). ---
$wsServices = $wb.Worksheets.Item("Services")
$wsServices.Range("A2").Select() | Out-Null

# --- FirewallPolicies sheet: collapse the old Protocol + Portnumbers
# columns into a single "Service" column that references the Services
# sheet, and rename the test row's Name. ---
$wsPolicies = $wb.Worksheets.Item("FirewallPolicies")
$wsPolicies.Range("G1").Value = "Service"
$wsPolicies.Range("G2").Value = "test_tcp"
$wsPolicies.Range("A2").Value = "test with ansible"

# Drop the now-unused Portnumbers column entirely (shifts dimension/
# spans from A1:H2 to A1:G2 automatically).
$wsPolicies.Columns("H:H").Delete()

# Make FirewallPolicies the active tab/sheet and park the selection on A4.
$wsPolicies.Activate() | Out-Null
$wsPolicies.Range("A4").Select() | Out-Null
